$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(125).Insert()
$ws.Range("A125").Value = 5
$ws.Range("B125").Value = "Macroferia Regional de Talca"
$ws.Range("C125").Value = "Maule"
$ws.Range("D125").Value = 45180
$ws.Range("E125").Value = 7
$ws.Range("F125").Value = 100112013
$ws.Range("G125").Value = "Alcachofa"
$ws.Range("H125").Value = "Madrigal"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 300
$ws.Range("K125").Value = 10000
$ws.Range("L125").Value = 10000
$ws.Range("M125").Value = 10000
$ws.Range("N125").Value = "$/caja 40 unidades"
$ws.Range("O125").Value = "Provincia del Elquí"
$ws.Range("P125").Value = 250
$ws.Range("Q125").Value = 40
$ws.Range("R125").Value = "Hortaliza"
